$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells
$ws.Range("G1").Value = "Total"
$ws.Range("H1").Value = "Attendance percentage"

# Mirror the bold/centered/top-aligned/bordered header style already used on A1:F1
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1:H1").HorizontalAlignment = -4108
$ws.Range("G1:H1").VerticalAlignment = -4160
$ws.Range("G1:H1").Borders.LineStyle = 1
$ws.Range("G1:H1").Borders.Weight = 2

# Per-row "Total" (count of Present across C:F) and "Attendance percentage"
$presentCounts = @{2 = 0; 3 = 1; 4 = 0; 5 = 0; 6 = 4}

foreach ($r in 2..6) {
    $total = $presentCounts[$r]
    $pct = $total / 4 * 100

    $ws.Cells.Item($r, 7).Value = $total
    $ws.Cells.Item($r, 8).Value = $pct
}
